# Turitea at No1 Dairy - river results refresh (May 2024)
# 1) Recalculated statistics for existing 2014-2018 / 2018-2022 rows (rows 3-96)
# 2) New 2019 - 2023 RepSite rows appended (rows 101-119)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: update recalculated values in existing rows ---
$cellUpdates = @(
    @{ Cell = "G3"; Value = 1.40731890673609 },
    @{ Cell = "I3"; Value = 3.0232 },
    @{ Cell = "L3"; Value = 1.6 },
    @{ Cell = "G4"; Value = 0.0208494286722883 },
    @{ Cell = "G5"; Value = 0.0208494286722883 },
    @{ Cell = "G6"; Value = 1876.9209492124 },
    @{ Cell = "H6"; Value = 11995.3966905906 },
    @{ Cell = "I6"; Value = 9890.5434 },
    @{ Cell = "G7"; Value = 1876.9209492124 },
    @{ Cell = "H7"; Value = 11995.3966905906 },
    @{ Cell = "I7"; Value = 9890.5434 },
    @{ Cell = "G8"; Value = 1876.9209492124 },
    @{ Cell = "H8"; Value = 11995.3966905906 },
    @{ Cell = "I8"; Value = 9890.5434 },
    @{ Cell = "G9"; Value = 1876.9209492124 },
    @{ Cell = "H9"; Value = 11995.3966905906 },
    @{ Cell = "I9"; Value = 9890.5434 },
    @{ Cell = "F11"; Value = 0.00171 },
    @{ Cell = "G11"; Value = 0.0063197838889798 },
    @{ Cell = "L11"; Value = 0.00186 },
    @{ Cell = "F12"; Value = 0.00171 },
    @{ Cell = "G12"; Value = 0.0063197838889798 },
    @{ Cell = "L12"; Value = 0.00186 },
    @{ Cell = "G13"; Value = 0.160716438133231 },
    @{ Cell = "L13"; Value = 0.00596 },
    @{ Cell = "G14"; Value = 0.160716438133231 },
    @{ Cell = "L14"; Value = 0.00596 },
    @{ Cell = "G16"; Value = 0.176442265325718 },
    @{ Cell = "H16"; Value = 0.6963 },
    @{ Cell = "N16"; Value = 0.48993 },
    @{ Cell = "G17"; Value = 0.176442265325718 },
    @{ Cell = "H17"; Value = 0.6963 },
    @{ Cell = "N17"; Value = 0.48993 },
    @{ Cell = "G23"; Value = 1.26109358644013 },
    @{ Cell = "I23"; Value = 2.617 },
    @{ Cell = "L23"; Value = 1.58034 },
    @{ Cell = "N23"; Value = 2.43838 },
    @{ Cell = "G24"; Value = 0.0212186048610834 },
    @{ Cell = "G25"; Value = 0.0212186048610834 },
    @{ Cell = "G26"; Value = 1771.14262803931 },
    @{ Cell = "H26"; Value = 11995.3966905906 },
    @{ Cell = "I26"; Value = 9872.83145 },
    @{ Cell = "G27"; Value = 1771.14262803931 },
    @{ Cell = "H27"; Value = 11995.3966905906 },
    @{ Cell = "I27"; Value = 9872.83145 },
    @{ Cell = "G28"; Value = 1771.14262803931 },
    @{ Cell = "H28"; Value = 11995.3966905906 },
    @{ Cell = "I28"; Value = 9872.83145 },
    @{ Cell = "G29"; Value = 1771.14262803931 },
    @{ Cell = "H29"; Value = 11995.3966905906 },
    @{ Cell = "I29"; Value = 9872.83145 },
    @{ Cell = "F31"; Value = 0.00215 },
    @{ Cell = "G31"; Value = 0.0074011415012237 },
    @{ Cell = "L31"; Value = 0.00205 },
    @{ Cell = "F32"; Value = 0.00215 },
    @{ Cell = "G32"; Value = 0.0074011415012237 },
    @{ Cell = "L32"; Value = 0.00205 },
    @{ Cell = "G33"; Value = 0.164462629793119 },
    @{ Cell = "G34"; Value = 0.164462629793119 },
    @{ Cell = "G36"; Value = 0.181273227570289 },
    @{ Cell = "N36"; Value = 0.48991 },
    @{ Cell = "G37"; Value = 0.181273227570289 },
    @{ Cell = "N37"; Value = 0.48991 },
    @{ Cell = "G43"; Value = 1.16758149399015 },
    @{ Cell = "H43"; Value = 2.94 },
    @{ Cell = "I43"; Value = 2.55702 },
    @{ Cell = "G44"; Value = 0.0207646829095764 },
    @{ Cell = "G45"; Value = 0.0207646829095764 },
    @{ Cell = "G46"; Value = 1436.5933630055 },
    @{ Cell = "H46"; Value = 11995.3966905906 },
    @{ Cell = "I46"; Value = 8099.3434 },
    @{ Cell = "G47"; Value = 1436.5933630055 },
    @{ Cell = "H47"; Value = 11995.3966905906 },
    @{ Cell = "I47"; Value = 8099.3434 },
    @{ Cell = "G48"; Value = 1436.5933630055 },
    @{ Cell = "H48"; Value = 11995.3966905906 },
    @{ Cell = "I48"; Value = 8099.3434 },
    @{ Cell = "G49"; Value = 1436.5933630055 },
    @{ Cell = "H49"; Value = 11995.3966905906 },
    @{ Cell = "I49"; Value = 8099.3434 },
    @{ Cell = "F51"; Value = 0.00279 },
    @{ Cell = "G51"; Value = 0.0065499333503473 },
    @{ Cell = "L51"; Value = 0.00298 },
    @{ Cell = "F52"; Value = 0.00279 },
    @{ Cell = "G52"; Value = 0.0065499333503473 },
    @{ Cell = "L52"; Value = 0.00298 },
    @{ Cell = "G53"; Value = 0.163465119707724 },
    @{ Cell = "G54"; Value = 0.163465119707724 },
    @{ Cell = "G56"; Value = 0.179518452277016 },
    @{ Cell = "G57"; Value = 0.179518452277016 },
    @{ Cell = "G63"; Value = 0.919041729883327 },
    @{ Cell = "H63"; Value = 2.53859344300821 },
    @{ Cell = "G66"; Value = 1634.00411248456 },
    @{ Cell = "H66"; Value = 11995.3966905906 },
    @{ Cell = "I66"; Value = 9179.3434 },
    @{ Cell = "G67"; Value = 1634.00411248456 },
    @{ Cell = "H67"; Value = 11995.3966905906 },
    @{ Cell = "I67"; Value = 9179.3434 },
    @{ Cell = "G68"; Value = 1634.00411248456 },
    @{ Cell = "H68"; Value = 11995.3966905906 },
    @{ Cell = "I68"; Value = 9179.3434 },
    @{ Cell = "G69"; Value = 1634.00411248456 },
    @{ Cell = "H69"; Value = 11995.3966905906 },
    @{ Cell = "I69"; Value = 9179.3434 },
    @{ Cell = "F71"; Value = 0.0046 },
    @{ Cell = "G71"; Value = 0.0074327040696234 },
    @{ Cell = "L71"; Value = 0.00438 },
    @{ Cell = "F72"; Value = 0.0046 },
    @{ Cell = "G72"; Value = 0.0074327040696234 },
    @{ Cell = "L72"; Value = 0.00438 },
    @{ Cell = "G73"; Value = 0.168245379184381 },
    @{ Cell = "G74"; Value = 0.168245379184381 },
    @{ Cell = "G76"; Value = 0.185304548913851 },
    @{ Cell = "G77"; Value = 0.185304548913851 },
    @{ Cell = "G85"; Value = 1726.29135960991 },
    @{ Cell = "H85"; Value = 11995.3966905906 },
    @{ Cell = "I85"; Value = 10148.48911 },
    @{ Cell = "G86"; Value = 1726.29135960991 },
    @{ Cell = "H86"; Value = 11995.3966905906 },
    @{ Cell = "I86"; Value = 10148.48911 },
    @{ Cell = "G87"; Value = 1726.29135960991 },
    @{ Cell = "H87"; Value = 11995.3966905906 },
    @{ Cell = "I87"; Value = 10148.48911 },
    @{ Cell = "G88"; Value = 1726.29135960991 },
    @{ Cell = "H88"; Value = 11995.3966905906 },
    @{ Cell = "I88"; Value = 10148.48911 },
    @{ Cell = "G90"; Value = 0.0085032470338379 },
    @{ Cell = "G91"; Value = 0.0085032470338379 },
    @{ Cell = "G92"; Value = 0.173160633421669 },
    @{ Cell = "G93"; Value = 0.173160633421669 },
    @{ Cell = "G95"; Value = 0.191507938744359 },
    @{ Cell = "G96"; Value = 0.191507938744359 }
)
foreach ($u in $cellUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Part 2: append new RepSite rows (2019 - 2023) as rows 101-119 ---
$newRows = @(
    @{ Row = 101; A="Turitea at No1 Dairy"; B="ASPM"; C="D"; D="2019 - 2023"; E="RepSite"; F=0.241; G=0.2882; H=0.45; I=0.45; J=$null; K=$null; L=0.241; M=0.415; N=0.45; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U=$null },
    @{ Row = 102; A="Turitea at No1 Dairy"; B="DRP (95th Percentile)"; C="C"; D="2019 - 2023"; E="RepSite"; F=0.017; G=0.0186271186440678; H=0.038; I=0.03655; J=$null; K=$null; L=0.023; M=0.02547; N=0.03378; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="mg/L" },
    @{ Row = 103; A="Turitea at No1 Dairy"; B="DRP (Median)"; C="C"; D="2019 - 2023"; E="RepSite"; F=0.017; G=0.0186271186440678; H=0.038; I=0.03655; J=$null; K=$null; L=0.023; M=0.02547; N=0.03378; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="mg/L" },
    @{ Row = 104; A="Turitea at No1 Dairy"; B="E coli (>260)"; C="E"; D="2019 - 2023"; E="RepSite"; F=450; G=1297.07719533314; H=10884.9516928798; I=7135; J=40.6779661016949; K=64.406779661017; L=490; M=1714.68; N=5554.12; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="% exceedances over 260/100 mL" },
    @{ Row = 105; A="Turitea at No1 Dairy"; B="E coli (>540)"; C="E"; D="2019 - 2023"; E="RepSite"; F=450; G=1297.07719533314; H=10884.9516928798; I=7135; J=40.6779661016949; K=64.406779661017; L=490; M=1714.68; N=5554.12; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="% exceedances over 540/100 mL" },
    @{ Row = 106; A="Turitea at No1 Dairy"; B="E coli (Median)"; C="E"; D="2019 - 2023"; E="RepSite"; F=450; G=1297.07719533314; H=10884.9516928798; I=7135; J=40.6779661016949; K=64.406779661017; L=490; M=1714.68; N=5554.12; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="E. coli/100 mL" },
    @{ Row = 107; A="Turitea at No1 Dairy"; B="E coli (95th Percentile)"; C="E"; D="2019 - 2023"; E="RepSite"; F=450; G=1297.07719533314; H=10884.9516928798; I=7135; J=40.6779661016949; K=64.406779661017; L=490; M=1714.68; N=5554.12; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="E. coli/100 mL" },
    @{ Row = 108; A="Turitea at No1 Dairy"; B="MCI"; C="C"; D="2019 - 2023"; E="RepSite"; F=92; G=92.834; H=109.17; I=109.17; J=$null; K=$null; L=92; M=105.2605; N=109.17; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U=$null },
    @{ Row = 109; A="Turitea at No1 Dairy"; B="Ammoniacal-N (95th Percentile)"; C="A"; D="2019 - 2023"; E="RepSite"; F=0.00486; G=0.0063895431772339; H=0.0372876372975684; I=0.0146; J=$null; K=$null; L=0.00596; M=0.00928; N=0.01328; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="mg NH4-N/L" },
    @{ Row = 110; A="Turitea at No1 Dairy"; B="Ammoniacal-N (Median)"; C="A"; D="2019 - 2023"; E="RepSite"; F=0.00486; G=0.0063895431772339; H=0.0372876372975684; I=0.0146; J=$null; K=$null; L=0.00596; M=0.00928; N=0.01328; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="mg NH4-N/L" },
    @{ Row = 111; A="Turitea at No1 Dairy"; B="Nitrate-N (95th Percentile)"; C="A"; D="2019 - 2023"; E="RepSite"; F=0.098; G=0.143658114144402; H=0.545; I=0.4398; J=$null; K=$null; L=0.025; M=0.33288; N=0.39568; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="mg NO3-N/L" },
    @{ Row = 112; A="Turitea at No1 Dairy"; B="Nitrate-N (Median)"; C="A"; D="2019 - 2023"; E="RepSite"; F=0.098; G=0.143658114144402; H=0.545; I=0.4398; J=$null; K=$null; L=0.025; M=0.33288; N=0.39568; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="mg NO3-N/L" },
    @{ Row = 113; A="Turitea at No1 Dairy"; B="QMCI"; C="D"; D="2019 - 2023"; E="RepSite"; F=4.35; G=4.5208; H=5.21; I=5.21; J=$null; K=$null; L=4.35; M=5.0385; N=5.21; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U=$null },
    @{ Row = 114; A="Turitea at No1 Dairy"; B="Soluble Inorganic Nitrogen (95th Percentile)"; C=$null; D="2019 - 2023"; E="RepSite"; F=0.115; G=0.158829972642664; H=0.574; I=0.4534; J=$null; K=$null; L=0.0425; M=0.34888; N=0.41994; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="g/m3" },
    @{ Row = 115; A="Turitea at No1 Dairy"; B="Soluble Inorganic Nitrogen (Median)"; C=$null; D="2019 - 2023"; E="RepSite"; F=0.115; G=0.158829972642664; H=0.574; I=0.4534; J=$null; K=$null; L=0.0425; M=0.34888; N=0.41994; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="g/m3" },
    @{ Row = 116; A="Turitea at No1 Dairy"; B="Total Nitrogen (95th Percentile)"; C=$null; D="2019 - 2023"; E="RepSite"; F=0.3; G=0.391864406779661; H=2.4; I=0.817; J=$null; K=$null; L=0.22; M=0.5141; N=0.719; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="g/m3" },
    @{ Row = 117; A="Turitea at No1 Dairy"; B="Total Nitrogen (Median)"; C=$null; D="2019 - 2023"; E="RepSite"; F=0.3; G=0.391864406779661; H=2.4; I=0.817; J=$null; K=$null; L=0.22; M=0.5141; N=0.719; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="g/m3" },
    @{ Row = 118; A="Turitea at No1 Dairy"; B="Total Phosphorus (95th Percentile)"; C=$null; D="2019 - 2023"; E="RepSite"; F=0.041; G=0.0481694915254237; H=0.184; I=0.11705; J=$null; K=$null; L=0.0475; M=0.05894; N=0.10076; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="g/m3" },
    @{ Row = 119; A="Turitea at No1 Dairy"; B="Total Phosphorus (Median)"; C=$null; D="2019 - 2023"; E="RepSite"; F=0.041; G=0.0481694915254237; H=0.184; I=0.11705; J=$null; K=$null; L=0.0475; M=0.05894; N=0.10076; O=1821276.44; P=5526193.49; Q="Palmerston North City"; R="Manawatū"; S="Lower Manawatu"; T="Mana_11b"; U="g/m3" }
)

foreach ($r in $newRows) {
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")) {
        $v = $r[$col]
        if ($null -ne $v) {
            $ws.Range("$col$($r.Row)").Value = $v
        }
    }
}
